$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44243
$ws.Range("I2").Value = "Especial"
$ws.Range("J2").Value = 300
$ws.Range("K2").Value = 12000
$ws.Range("L2").Value = 12000
$ws.Range("M2").Value = 12000
$ws.Range("P2").Value = 667

# Row 3
$ws.Range("D3").Value = 44243
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 300
$ws.Range("K3").Value = 10000
$ws.Range("L3").Value = 10000
$ws.Range("M3").Value = 10000
$ws.Range("P3").Value = 556

# Row 4
$ws.Range("D4").Value = 44243
$ws.Range("I4").Value = "Segunda"
$ws.Range("J4").Value = 150
$ws.Range("K4").Value = 8000
$ws.Range("L4").Value = 8000
$ws.Range("M4").Value = 8000
$ws.Range("P4").Value = 444

# Row 5
$ws.Range("D5").Value = 44383
$ws.Range("J5").Value = 300
$ws.Range("K5").Value = 16000
$ws.Range("L5").Value = 16000
$ws.Range("M5").Value = 16000
$ws.Range("P5").Value = 889

# Row 6
$ws.Range("D6").Value = 44383
$ws.Range("K6").Value = 12000
$ws.Range("L6").Value = 12000
$ws.Range("M6").Value = 12000
$ws.Range("P6").Value = 667

# Row 8
$ws.Range("D8").Value = 44245
$ws.Range("J8").Value = 300
$ws.Range("K8").Value = 12000
$ws.Range("L8").Value = 12000
$ws.Range("M8").Value = 12000
$ws.Range("P8").Value = 667

# Row 9
$ws.Range("D9").Value = 44245
$ws.Range("J9").Value = 200
$ws.Range("K9").Value = 10000
$ws.Range("L9").Value = 10000
$ws.Range("M9").Value = 10000
$ws.Range("P9").Value = 556

# Row 10
$ws.Range("D10").Value = 44396
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 250
$ws.Range("K10").Value = 15000
$ws.Range("L10").Value = 15000
$ws.Range("M10").Value = 15000
$ws.Range("P10").Value = 833

# Row 11
$ws.Range("D11").Value = 44396
$ws.Range("I11").Value = "Segunda"
$ws.Range("J11").Value = 150
$ws.Range("K11").Value = 12000
$ws.Range("L11").Value = 12000
$ws.Range("M11").Value = 12000
$ws.Range("P11").Value = 667

# Row 12
$ws.Range("D12").Value = 44249
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 400
$ws.Range("K12").Value = 12000
$ws.Range("L12").Value = 12000
$ws.Range("M12").Value = 12000
$ws.Range("P12").Value = 667

# Row 13
$ws.Range("D13").Value = 44249
$ws.Range("I13").Value = "Segunda"
$ws.Range("J13").Value = 200
$ws.Range("K13").Value = 10000
$ws.Range("L13").Value = 10000
$ws.Range("M13").Value = 10000
$ws.Range("P13").Value = 556

# Row 14
$ws.Range("D14").Value = 44235
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 400
$ws.Range("K14").Value = 13000
$ws.Range("L14").Value = 13000
$ws.Range("M14").Value = 13000
$ws.Range("P14").Value = 722

# Row 15
$ws.Range("D15").Value = 44235
$ws.Range("I15").Value = "Segunda"
$ws.Range("J15").Value = 200
$ws.Range("K15").Value = 11000
$ws.Range("L15").Value = 11000
$ws.Range("M15").Value = 11000
$ws.Range("P15").Value = 611

# Row 16
$ws.Range("D16").Value = 44235
$ws.Range("I16").Value = "Tercera"
$ws.Range("J16").Value = 100
$ws.Range("K16").Value = 9000
$ws.Range("L16").Value = 9000
$ws.Range("M16").Value = 9000
$ws.Range("P16").Value = 500

# Row 17
$ws.Range("D17").Value = 44391
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 400
$ws.Range("K17").Value = 15000
$ws.Range("L17").Value = 15000
$ws.Range("M17").Value = 15000
$ws.Range("P17").Value = 833

# Row 18
$ws.Range("D18").Value = 44238
$ws.Range("J18").Value = 300
$ws.Range("K18").Value = 12000
$ws.Range("L18").Value = 12000
$ws.Range("M18").Value = 12000
$ws.Range("P18").Value = 667

# Row 19
$ws.Range("D19").Value = 44238
$ws.Range("I19").Value = "Segunda"
$ws.Range("J19").Value = 200
$ws.Range("K19").Value = 10000
$ws.Range("L19").Value = 10000
$ws.Range("M19").Value = 10000
$ws.Range("P19").Value = 556

# Row 20
$ws.Range("D20").Value = 44238
$ws.Range("I20").Value = "Tercera"
$ws.Range("J20").Value = 50
$ws.Range("K20").Value = 8000
$ws.Range("L20").Value = 8000
$ws.Range("M20").Value = 8000
$ws.Range("P20").Value = 444
